# Fix the mangled start-date in the last schedule row (2025/11/35 is not a
# valid date, it should be 2025/12/5) and add the period-87 description that
# was missing from column C for that row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A37").Value = "2025/12/5"
$ws.Range("C37").Value = "第87期 秘寶 開放區域 玫瑰淚堡 祕寶效果: 貓貓包融合成功5%(18%)增加一個招牌貓"
